# Applies the authoring changes captured in the diff:
#  - "Excel Notes" sheet gets a new value in C2 (=1)
#  - Selection on "Excel Notes" moves from C2 to B3
#  - The active/selected tab moves from "Excel Notes" to "ABC Notes"
#  - Window geometry is updated to reflect the new window placement

$wb = $excel.ActiveWorkbook

$wsExcelNotes = $wb.Worksheets.Item("Excel Notes")
$wsAbcNotes   = $wb.Worksheets.Item("ABC Notes")

# Work on the "Excel Notes" sheet first: add the new data point and move
# the selection there (this sheet is active while we do this).
$wsExcelNotes.Activate()
$wsExcelNotes.Range("C2").Value = 1
$wsExcelNotes.Range("B3").Select()

# Resize / reposition the workbook window to match the saved view state.
$aw = $excel.ActiveWindow
$aw.Left = -23148
$aw.Top = 9684
$aw.Width = 23256
$aw.Height = 13896

# Finally, make "ABC Notes" the active/selected sheet tab.
$wsAbcNotes.Activate()
